$d = $word.ActiveDocument

# --- 1. Title heading: "SCaR_Arcade" -> "SCaR Arcade" (scope to first paragraph only,
#        so the unrelated "SCaR_Arcade" mid-document is left alone) ---
$d.Paragraphs(1).Range.Find.Execute("SCaR_Arcade", $true, $false, $false, $false, $false, $true, 1, $false, "SCaR Arcade", 2)

# --- 2. "What is SCaR_Arcade?" -> "What is SCaR Arcade?" ---
$d.Content.Find.Execute("What is SCaR_Arcade?", $true, $false, $false, $false, $false, $true, 1, $false, "What is SCaR Arcade?", 2)

# --- 3. Intro paragraph: "prototype for a game application manger designed for the Android OS"
#        -> "Library Game Application Manager (LGAM) designed, and built for the Android OS"
#        (the leading "SCaR_Arcade" earlier in this same paragraph is intentionally left as-is) ---
$d.Content.Find.Execute("prototype for a game application manger designed for the Android OS", $true, $false, $false, $false, $false, $true, 1, $false, "Library Game Application Manager (LGAM) designed, and built for the Android OS", 2)

# --- 4. "Features of SCaR_Arcade" -> "Features of SCaR Arcade" ---
$d.Content.Find.Execute("Features of SCaR_Arcade", $true, $false, $false, $false, $false, $true, 1, $false, "Features of SCaR Arcade", 2)

# --- 5. "for each game application." -> "for each game." ---
$d.Content.Find.Execute("for each game application.", $true, $false, $false, $false, $false, $true, 1, $false, "for each game.", 2)

# --- 6. Remove the whole "Compares other players score globally..." bullet paragraph ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -like "Compares other players score globally*") {
        $para.Range.Delete()
        break
    }
}

# --- 7. "The main menu of SCaR_Arcade" -> "The main menu of SCaR Arcade" ---
$d.Content.Find.Execute("The main menu of SCaR_Arcade", $true, $false, $false, $false, $false, $true, 1, $false, "The main menu of SCaR Arcade", 2)

# --- 8. "...users choice to the game menu." -> "...users choice to the game menu activity." ---
$d.Content.Find.Execute("users choice to the game menu.", $true, $false, $false, $false, $false, $true, 1, $false, "users choice to the game menu activity.", 2)

# --- 9. Move the _GoBack bookmark to sit right before the final "." of "game menu activity." ---
$r = $d.Content
$r.Find.Execute("game menu activity", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r)

# --- 10. Leaderboard paragraph rewrite ---
$d.Content.Find.Execute("The Leaderboard Hub", $true, $false, $false, $false, $false, $true, 1, $false, "The leaderboard hub", 2)
$d.Content.Find.Execute("local high scores for the game and", $true, $false, $false, $false, $false, $true, 1, $false, "local high scores for a game and", 2)
$d.Content.Find.Execute("view the Online leaderboard", $true, $false, $false, $false, $false, $true, 1, $false, "view the online leaderboard", 2)

Write-Output "done"
